$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-21 (row index, then 20 values for columns A..T)
$data = @(
  ,@(2,"ECs","Ccl11","Ccr5","ECs",[double]"3",[double]"1",[double]"1.819052333333333",[double]"5.457157",[double]"0.02872704074300508",[double]"0.02872704074300508",[double]"3",[double]"1",[double]"0.1207436666666667",[double]"0.362231",[double]"0.001088347656299082",[double]"0.001088347656299082",[double]"0.2196390485852222",[double]"1.976751437267",[double]"3.126500746505783E-05",[double]"3.126500746505783E-05")
  ,@(3,"ECs","Ccl11","Ccr5","M1",[double]"3",[double]"1",[double]"1.819052333333333",[double]"5.457157",[double]"0.02872704074300508",[double]"0.02872704074300508",[double]"3",[double]"1",[double]"48.86031499999999",[double]"146.580945",[double]"0.4404124107236948",[double]"0.4404124107236948",[double]"88.87947000815166",[double]"799.9152300733649",[double]"0.01265174526658467",[double]"0.01265174526658467")
  ,@(4,"ECs","Ccl11","Ccr5","M2",[double]"3",[double]"1",[double]"1.819052333333333",[double]"5.457157",[double]"0.02872704074300508",[double]"0.02872704074300508",[double]"3",[double]"1",[double]"61.90825266666667",[double]"185.724758",[double]"0.5580226570503747",[double]"0.5580226570503747",[double]"112.6143514658896",[double]"1013.529163193006",[double]"0.01603033960460606",[double]"0.01603033960460606")
  ,@(5,"ECs","Ccl11","Ccr5","sCs",[double]"3",[double]"1",[double]"1.819052333333333",[double]"5.457157",[double]"0.02872704074300508",[double]"0.02872704074300508",[double]"2",[double]"0.6666666666666666",[double]"0.05287333333333333",[double]"0.15862",[double]"0.0004765845696314243",[double]"0.0004765845696314243",[double]"0.09617936037111111",[double]"0.86561424334",[double]"1.369086434928947E-05",[double]"1.369086434928947E-05")
  ,@(6,"FAPs","Ccl11","Ccr5","ECs",[double]"3",[double]"1",[double]"56.85979966666667",[double]"170.579399",[double]"0.8979476575422553",[double]"0.8979476575422553",[double]"3",[double]"1",[double]"0.1207436666666667",[double]"0.362231",[double]"0.001088347656299082",[double]"0.001088347656299082",[double]"6.865460697685444",[double]"61.789146279169",[double]"0.0009772792285653645",[double]"0.0009772792285653645")
  ,@(7,"FAPs","Ccl11","Ccr5","M1",[double]"3",[double]"1",[double]"56.85979966666667",[double]"170.579399",[double]"0.8979476575422553",[double]"0.8979476575422553",[double]"3",[double]"1",[double]"48.86031499999999",[double]"146.580945",[double]"0.4404124107236948",[double]"0.4404124107236948",[double]"2778.187722550229",[double]"25003.68950295206",[double]"0.3954672925618793",[double]"0.3954672925618793")
  ,@(8,"FAPs","Ccl11","Ccr5","M2",[double]"3",[double]"1",[double]"56.85979966666667",[double]"170.579399",[double]"0.8979476575422553",[double]"0.8979476575422553",[double]"3",[double]"1",[double]"61.90825266666667",[double]"185.724758",[double]"0.5580226570503747",[double]"0.5580226570503747",[double]"3520.09084434005",[double]"31680.81759906045",[double]"0.5010751377538892",[double]"0.5010751377538892")
  ,@(9,"FAPs","Ccl11","Ccr5","sCs",[double]"3",[double]"1",[double]"56.85979966666667",[double]"170.579399",[double]"0.8979476575422553",[double]"0.8979476575422553",[double]"2",[double]"0.6666666666666666",[double]"0.05287333333333333",[double]"0.15862",[double]"0.0004765845696314243",[double]"0.0004765845696314243",[double]"3.006367141042223",[double]"27.05730426938",[double]"0.0004279479979213213",[double]"0.0004279479979213213")
  ,@(10,"M1","Ccl11","Ccr5","ECs",[double]"3",[double]"1",[double]"1.611504",[double]"4.834512",[double]"0.02544937285046902",[double]"0.02544937285046902",[double]"3",[double]"1",[double]"0.1207436666666667",[double]"0.362231",[double]"0.001088347656299082",[double]"0.001088347656299082",[double]"0.194578901808",[double]"1.751210116272",[double]"2.769776529608946E-05",[double]"2.769776529608946E-05")
  ,@(11,"M1","Ccl11","Ccr5","M1",[double]"3",[double]"1",[double]"1.611504",[double]"4.834512",[double]"0.02544937285046902",[double]"0.02544937285046902",[double]"3",[double]"1",[double]"48.86031499999999",[double]"146.580945",[double]"0.4404124107236948",[double]"0.4404124107236948",[double]"78.73859306375999",[double]"708.64733757384",[double]"0.01120821964848121",[double]"0.01120821964848121")
  ,@(12,"M1","Ccl11","Ccr5","M2",[double]"3",[double]"1",[double]"1.611504",[double]"4.834512",[double]"0.02544937285046902",[double]"0.02544937285046902",[double]"3",[double]"1",[double]"61.90825266666667",[double]"185.724758",[double]"0.5580226570503747",[double]"0.5580226570503747",[double]"99.76539680534401",[double]"897.888571248096",[double]"0.01420132665828439",[double]"0.01420132665828439")
  ,@(13,"M1","Ccl11","Ccr5","sCs",[double]"3",[double]"1",[double]"1.611504",[double]"4.834512",[double]"0.02544937285046902",[double]"0.02544937285046902",[double]"2",[double]"0.6666666666666666",[double]"0.05287333333333333",[double]"0.15862",[double]"0.0004765845696314243",[double]"0.0004765845696314243",[double]"0.08520558816",[double]"0.7668502934400001",[double]"1.212877840733043E-05",[double]"1.212877840733043E-05")
  ,@(14,"M2","Ccl11","Ccr5","ECs",[double]"3",[double]"1",[double]"1.78513",[double]"5.35539",[double]"0.0281913286945349",[double]"0.0281913286945349",[double]"3",[double]"1",[double]"0.1207436666666667",[double]"0.362231",[double]"0.001088347656299082",[double]"0.001088347656299082",[double]"0.2155431416766666",[double]"1.93988827509",[double]"3.068196651265412E-05",[double]"3.068196651265412E-05")
  ,@(15,"M2","Ccl11","Ccr5","M1",[double]"3",[double]"1",[double]"1.78513",[double]"5.35539",[double]"0.0281913286945349",[double]"0.0281913286945349",[double]"3",[double]"1",[double]"48.86031499999999",[double]"146.580945",[double]"0.4404124107236948",[double]"0.4404124107236948",[double]"87.22201411594997",[double]"784.9981270435499",[double]"0.01241581103186419",[double]"0.01241581103186419")
  ,@(16,"M2","Ccl11","Ccr5","M2",[double]"3",[double]"1",[double]"1.78513",[double]"5.35539",[double]"0.0281913286945349",[double]"0.0281913286945349",[double]"3",[double]"1",[double]"61.90825266666667",[double]"185.724758",[double]"0.5580226570503747",[double]"0.5580226570503747",[double]"110.5142790828467",[double]"994.62851174562",[double]"0.01573140014390483",[double]"0.01573140014390484")
  ,@(17,"M2","Ccl11","Ccr5","sCs",[double]"3",[double]"1",[double]"1.78513",[double]"5.35539",[double]"0.0281913286945349",[double]"0.0281913286945349",[double]"2",[double]"0.6666666666666666",[double]"0.05287333333333333",[double]"0.15862",[double]"0.0004765845696314243",[double]"0.0004765845696314243",[double]"0.09438577353333333",[double]"0.8494719618000001",[double]"1.343555225322294E-05",[double]"1.343555225322294E-05")
  ,@(18,"sCs","Ccl11","Ccr5","ECs",[double]"3",[double]"1",[double]"1.246467333333333",[double]"3.739402",[double]"0.01968460016973576",[double]"0.01968460016973576",[double]"3",[double]"1",[double]"0.1207436666666667",[double]"0.362231",[double]"0.001088347656299082",[double]"0.001088347656299082",[double]"0.1505030362068889",[double]"1.354527325862",[double]"2.142368845991643E-05",[double]"2.142368845991643E-05")
  ,@(19,"sCs","Ccl11","Ccr5","M1",[double]"3",[double]"1",[double]"1.246467333333333",[double]"3.739402",[double]"0.01968460016973576",[double]"0.01968460016973576",[double]"3",[double]"1",[double]"48.86031499999999",[double]"146.580945",[double]"0.4404124107236948",[double]"0.4404124107236948",[double]"60.90278654387666",[double]"548.12507889489",[double]"0.008669342214885379",[double]"0.008669342214885379")
  ,@(20,"sCs","Ccl11","Ccr5","M2",[double]"3",[double]"1",[double]"1.246467333333333",[double]"3.739402",[double]"0.01968460016973576",[double]"0.01968460016973576",[double]"3",[double]"1",[double]"61.90825266666667",[double]"185.724758",[double]"0.5580226570503747",[double]"0.5580226570503747",[double]"77.16661461274623",[double]"694.4995315147161",[double]"0.01098445288969021",[double]"0.01098445288969021")
  ,@(21,"sCs","Ccl11","Ccr5","sCs",[double]"3",[double]"1",[double]"1.246467333333333",[double]"3.739402",[double]"0.01968460016973576",[double]"0.01968460016973576",[double]"2",[double]"0.6666666666666666",[double]"0.05287333333333333",[double]"0.15862",[double]"0.0004765845696314243",[double]"0.0004765845696314243",[double]"0.06590488280444445",[double]"0.59314394524",[double]"9.38137670026018E-06",[double]"9.38137670026018E-06")
)
foreach ($rowEntry in $data) {
    $rowNum = $rowEntry[0]
    for ($colIdx = 1; $colIdx -le 20; $colIdx++) {
        $ws.Cells.Item($rowNum, $colIdx).Value = $rowEntry[$colIdx]
    }
}
